# Fruta / hortaliza, semanal
# The weekly refresh reshuffles which historical record occupies each row
# (columns D, L, M, N, O, P, Q, R, S, T) while the descriptive columns
# (A, B, C, E, F, G, H, I, J, K) stay put. Row 2..41 -> source row mapping:

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2 = 9;  3 = 10; 4 = 12; 5 = 5;  6 = 28; 7 = 20; 8 = 22; 9 = 6;  10 = 13;
    11 = 27; 12 = 8; 13 = 37; 14 = 39; 15 = 40; 16 = 31; 17 = 2; 18 = 23;
    19 = 25; 20 = 16; 21 = 24; 22 = 15; 23 = 17; 24 = 18; 25 = 19; 26 = 29;
    27 = 30; 28 = 14; 29 = 7; 30 = 36; 31 = 41; 32 = 11; 33 = 38; 34 = 32;
    35 = 33; 36 = 26; 37 = 3; 38 = 4; 39 = 34; 40 = 21; 41 = 35
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the moving columns for every data row BEFORE writing anything,
# since several destinations are also sources for other rows.
$snapshot = @{}
foreach ($r in 2..41) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in 2..41) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
